$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same bold/border/alignment style (s="1").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data row values
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8
